$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A, shifting existing data right by 2 columns
$ws.Range("A:B").Insert()

# Fill new column A with reference, column B with table, for data rows
$lastRow = 17

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "Barlow et al. 1994"
    $ws.Cells.Item($r, 2).Value = "Table 3"
}

$ws.Range("A1").Value = "reference"
$ws.Range("B1").Value = "table"

# Apply bold style (same as header cells) to A1:B1 to match header row style
$ws.Range("A1:B1").Font.Bold = $true

# Autofit new columns A and B
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

# Update selection to match target state
$ws.Range("A15").Select()
